$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Payments")
$pt = $ws.PivotTables(1)
$pf = $pt.PivotFields("Payment Date")
Write-Host "Orientation before:" $pf.Orientation
$items = $pf.PivotItems()
Write-Host "Count:" $items.Count
# Try ShowDetail approach while staying in row orientation
$it1 = $items.Item(1)
try {
  $it1.ShowDetail = $false
  Write-Host "showdetail set ok"
} catch { Write-Host "err:" $_ }
